# The commit swaps the contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: theme2.xml (the theme actually used by the slide
# master / presentation) goes from the "Integral" design's "Red Violet"
# colour scheme back to the plain default "Office" colour scheme, while
# theme1.xml (only ever linked from the notes master) picks up the
# "Integral"/"Red Violet" colours that theme2.xml used to hold. The
# <a:fontScheme>/<a:fmtScheme> blocks of both theme parts are already
# byte-for-byte identical, so the only observable content difference is
# the 12-slot colour scheme (plus the cosmetic name="" attributes, which
# PowerPoint's object model does not expose a way to rewrite).
#
# Re-apply that by pushing the "Office" colour values onto the theme that
# the presentation's Design/SlideMaster actually exposes.

$p = $ppt.ActivePresentation

$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Order matches the DrawingML <a:clrScheme> child order: dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink. Values are packed as COM RGB
# (0x00BBGGRR = R + G*256 + B*65536) straight from the target theme's hex
# srgbClr values.
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
